$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "56.196.55"
$ws.Range("D3").Value = "2.982.81"
$ws.Range("E3").Value = "  -0.27%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").Value = "504.68"
$ws.Range("E5").Value = "  +0.87%  "
$ws.Range("D6").Value = "'137.20"
$ws.Range("E6").Value = "  -0.83%  "
$ws.Range("E7").Value = "  +0.11%  "
$ws.Range("D8").Value = "0.427"
$ws.Range("E8").Value = "  -0.78%  "
$ws.Range("E9").Value = "  -1.71%  "
$ws.Range("E10").Value = "  -1.52%  "
$ws.Range("D11").Value = "0.365"
$ws.Range("E11").Value = "  +1.15%  "
$ws.Range("D12").Value = "3.488.53"
$ws.Range("E12").Value = "  -0.21%  "
$ws.Range("E13").Value = "  -1.25%  "
$ws.Range("D14").Value = "25.89"
$ws.Range("E14").Value = "  -1.21%  "
$ws.Range("E15").Value = "  +0.47%  "
$ws.Range("D16").Value = "56.119.28"
$ws.Range("E16").Value = "  -1.64%  "
$ws.Range("D17").Value = "2.975.98"
$ws.Range("E17").Value = "  -0.40%  "
$ws.Range("D18").Value = "5.98"
$ws.Range("E18").Value = "  -2.07%  "
$ws.Range("D19").Value = "12.85"
$ws.Range("E19").Value = "  +1.28%  "
$ws.Range("D20").Value = "7.97"
$ws.Range("E20").Value = "  +0.74%  "
$ws.Range("D21").Value = "328.89"
$ws.Range("E21").Value = "  +1.77%  "
$ws.Range("E22").Value = "  +0.18%  "
$ws.Range("D23").Value = "0.493"
$ws.Range("E23").Value = "  -0.03%  "
$ws.Range("D24").Value = "64.44"
$ws.Range("E24").Value = "  +1.46%  "
$ws.Range("D25").Value = "3.101.24"
$ws.Range("D26").Value = "'1.00"
$ws.Range("E26").Value = "  -0.17%  "
$ws.Range("E27").Value = "  -0.71%  "
$ws.Range("D28").Value = "0.0₃0918"
$ws.Range("E28").Value = "  +2.33%  "
$ws.Range("D29").Value = "6.35"
$ws.Range("E29").Value = "  -3.62%  "
$ws.Range("E30").Value = "  -1.71%  "
$ws.Range("D31").Value = "1.78"
$ws.Range("E31").Value = "  +0.45%  "
$ws.Range("D32").Value = "1.16"
$ws.Range("E32").Value = "  -1.09%  "
$ws.Range("D33").Value = "20.11"
$ws.Range("E33").Value = "  -1.14%  "
$ws.Range("D34").Value = "152.72"
$ws.Range("E34").Value = "  -1.92%  "
$ws.Range("D35").Value = "4.49"
$ws.Range("E35").Value = "  -2.24%  "
$ws.Range("D36").Value = "5.77"
$ws.Range("E36").Value = "  -0.85%  "
$ws.Range("D37").Value = "25.69"
$ws.Range("E37").Value = "  +5.42%  "
$ws.Range("D38").Value = "1.24"
$ws.Range("E38").Value = "  -1.58%  "
$ws.Range("D39").Value = "0.0659"
$ws.Range("E39").Value = "  -1.40%  "
$ws.Range("D40").Value = "3.013.18"
$ws.Range("E40").Value = "  -0.09%  "
$ws.Range("D41").Value = "36.88"
$ws.Range("E41").Value = "  -2.49%  "
$ws.Range("D42").Value = "0.999"
$ws.Range("E42").Value = "  -0.07%  "
$ws.Range("D43").Value = "3.78"
$ws.Range("E43").Value = "  +0.36%  "
$ws.Range("E44").Value = "  +0.65%  "
$ws.Range("D45").Value = "2.169.58"
$ws.Range("E45").Value = "  -1.75%  "
$ws.Range("E46").Value = "  -3.02%  "
$ws.Range("D47").Value = "5.82"
$ws.Range("E47").Value = "  -2.80%  "
$ws.Range("B48").Value = "ONDO"
$ws.Range("C48").Value = "https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo"
$ws.Range("D48").Value = "0.919"
$ws.Range("E48").Value = "  -3.05%  "
$ws.Range("B49").Value = "VeChain"
$ws.Range("C49").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D49").Value = "0.0236"
$ws.Range("E49").Value = "  +0.37%  "
$ws.Range("D50").Value = "19.43"
$ws.Range("E50").Value = "  +0.26%  "
$ws.Range("D51").Value = "0.0847"
$ws.Range("E51").Value = "  -3.45%  "
